# Added Date Filter, Hid Cluster8, improved empty plot logic
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the descriptor2/descriptor3 columns for the "Habilidad Comercial" row
# (cluster 8, row 21) - previously filled descriptor cells now left empty.
$ws.Range("E21:G21").ClearContents()

# Update the view: scroll so column C is the left-most visible column, and
# move the active selection to C21.
$ws.Range("C21").Select()
$excel.ActiveWindow.ScrollColumn = 3
